# dati_quarra.xlsx - weekly update
# - Produzione: correct the 45779 (week of 28/07) production figure
# - Entrate / Uscite / Saldo: append the new week (45786) row of figures
# - leave the user's cursor / view state where it ended up in each sheet

$wb = $excel.ActiveWorkbook

# --- Produzione: fix row 19 (date 45779) value ---
$wsP = $wb.Worksheets.Item("Produzione")
$wsP.Range("B19").Value = 19124.98

# --- Entrate: new row 20 (date 45786) ---
$wsE = $wb.Worksheets.Item("Entrate")
$wsE.Range("A20").Value = 45786
$wsE.Range("B20").Value = 14800
$wsE.Range("B20").NumberFormat = "0.00"

# --- Uscite: new row 20 (date 45786) ---
$wsU = $wb.Worksheets.Item("Uscite")
$wsU.Range("A20").Value = 45786
$wsU.Range("B20").Value = 31875.47
$wsU.Range("B20").NumberFormat = "0.00"

# --- Saldo: new row 20 (date 45786) ---
$wsS = $wb.Worksheets.Item("Saldo")
$wsS.Range("A20").Value = 45786
$wsS.Range("B20").Value = 39070

# --- restore each sheet's on-screen selection, in the order the author last
#     touched them, so the final "tabSelected"/ActiveWindow state matches ---
$wsE.Range("A20:B20").Select()

$wsU.Range("E34").Select()

$wsS.Activate()
$excel.ActiveWindow.Zoom = 110
$wsS.Range("A20:B20").Select()

$wsP.Activate()
$wsP.Range("B20").Select()
